$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 500.84375
$ws.Range("I4").Value = 517.36664
$ws.Range("J4").Value = 253
$ws.Range("K4").Value = 517.36664
$ws.Range("L4").Value = 253
$ws.Range("M4").Value = -403.36664
$ws.Range("N4").Value = -481

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 635.1429000000001
$ws.Range("I28").Value = 569.4
$ws.Range("J28").Value = 799.5
$ws.Range("K28").Value = 569.4
$ws.Range("L28").Value = 799.5
$ws.Range("M28").Value = -84.39999999999998
$ws.Range("N28").Value = -1769.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 499.66666
$ws.Range("I41").Value = 499.66666
$ws.Range("K41").Value = 499.66666
$ws.Range("M41").Value = -59.66665999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 749.8570999999999
$ws.Range("I53").Value = 739.8
$ws.Range("K53").Value = 739.8
$ws.Range("M53").Value = -102.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3766.3333
$ws.Range("I62").Value = 2649.5
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 2649.5
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -2025.5
$ws.Range("N62").Value = -7248

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3766.3333
$ws.Range("I65").Value = 2649.5
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 13247.5
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -10127.5
$ws.Range("N65").Value = -36240

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5997.1665
$ws.Range("I76").Value = 5997.1665
$ws.Range("K76").Value = 5997.1665
$ws.Range("M76").Value = -5682.1665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5997.1665
$ws.Range("I79").Value = 5997.1665
$ws.Range("K79").Value = 5997.1665
$ws.Range("M79").Value = -4905.1665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 2099.8
$ws.Range("I92").Value = 2099.8
$ws.Range("K92").Value = 2099.8
$ws.Range("M92").Value = -851.8000000000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5529
$ws.Range("I98").Value = 1398.5714
$ws.Range("J98").Value = 15166.667
$ws.Range("K98").Value = 1398.5714
$ws.Range("L98").Value = 15166.667
$ws.Range("M98").Value = 99.42859999999996
$ws.Range("N98").Value = -18162.667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3000
$ws.Range("I106").Value = 3000
$ws.Range("K106").Value = 3000
$ws.Range("M106").Value = -2369

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 466.33334
$ws.Range("I107").Value = 199.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 199.5
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1720.5
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2723.4
$ws.Range("J112").Value = 2723.4
$ws.Range("L112").Value = 8170.200000000001
$ws.Range("N112").Value = -10386.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 5529
$ws.Range("I122").Value = 1398.5714
$ws.Range("J122").Value = 15166.667
$ws.Range("K122").Value = 4195.7142
$ws.Range("L122").Value = 45500.001
$ws.Range("M122").Value = -1745.7142
$ws.Range("N122").Value = -50400.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10858.857
$ws.Range("I32").Value = 7335.3335
$ws.Range("K32").Value = 7335.3335
$ws.Range("M32").Value = -7048.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1679.5714
$ws.Range("I132").Value = 1404
$ws.Range("K132").Value = 4212
$ws.Range("M132").Value = -1682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2040.1
$ws.Range("I134").Value = 1526.1111
$ws.Range("K134").Value = 4578.3333
$ws.Range("M134").Value = -2043.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1080.9412
$ws.Range("I2").Value = 418.33334
$ws.Range("J2").Value = 6050.5
$ws.Range("K2").Value = 2510.00004
$ws.Range("L2").Value = 36303
$ws.Range("M2").Value = -2397.00004
$ws.Range("N2").Value = -36529

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 675
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 675
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 2025
$ws.Range("N81").Value = -4271
$ws.Range("M81").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 675
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 675
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 6075
$ws.Range("N84").Value = -17307
$ws.Range("M84").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1633.4
$ws.Range("I140").Value = 1633.4
$ws.Range("K140").Value = 4900.200000000001
$ws.Range("M140").Value = 279.7999999999993

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.85714
$ws.Range("I2").Value = 59.666668
$ws.Range("K2").Value = 59.666668
$ws.Range("M2").Value = 53.333332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5760
$ws.Range("I5").Value = 40
$ws.Range("J5").Value = 7666.6665
$ws.Range("K5").Value = 40
$ws.Range("L5").Value = 7666.6665
$ws.Range("M5").Value = 72
$ws.Range("N5").Value = -7890.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 75
$ws.Range("I46").Value = 75
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 75
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 81
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 4999
$ws.Range("I5").Value = 4999
$ws.Range("K5").Value = 4999
$ws.Range("M5").Value = -4886

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2803.8
$ws.Range("J122").Value = 3875
$ws.Range("L122").Value = 11625
$ws.Range("N122").Value = -16525

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2021.3334
$ws.Range("I113").Value = 1058.8334
$ws.Range("J113").Value = 2983.8333
$ws.Range("K113").Value = 3176.5002
$ws.Range("L113").Value = 8951.499899999999
$ws.Range("M113").Value = -1006.5002
$ws.Range("N113").Value = -13291.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4591
$ws.Range("I122").Value = 4350
$ws.Range("K122").Value = 13050
$ws.Range("M122").Value = -10600
